$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "43.020.30"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.310.68"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "304.21"
$ws.Range("E5").Value = "  +2.07%  "
Set-TextValue $ws.Range("D6") "100.46"
$ws.Range("E6").Value = "  +5.44%  "
$ws.Range("E7").Value = "  +2.79%  "
Set-TextValue $ws.Range("D9") "0.513"
$ws.Range("E9").Value = "  +4.20%  "
Set-TextValue $ws.Range("D10") "34.94"
$ws.Range("E10").Value = "  +4.44%  "
Set-TextValue $ws.Range("D11") "0.0797"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("E12").Value = "  +4.19%  "
Set-TextValue $ws.Range("D13") "18.13"
$ws.Range("E13").Value = "  +16.71%  "
Set-TextValue $ws.Range("D14") "6.89"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").Value = "2.687.95"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").Value = "2.345.13"
$ws.Range("E16").Value = "  +3.11%  "
Set-TextValue $ws.Range("D17") "0.817"
$ws.Range("E17").Value = "  +4.44%  "
$ws.Range("D18").Value = "42.962.02"
$ws.Range("E18").Value = "  +1.99%  "
Set-TextValue $ws.Range("D19") "12.52"
$ws.Range("E19").Value = "  +6.99%  "
Set-TextValue $ws.Range("D20") "6.17"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("E22").Value = "  +1.92%  "
Set-TextValue $ws.Range("D23") "237.31"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("E24").Value = "  +12.75%  "
Set-TextValue $ws.Range("D25") "2.47"
$ws.Range("E25").Value = "  +0.70%  "
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  -0.43%  "
Set-TextValue $ws.Range("D27") "24.85"
$ws.Range("E27").Value = "  +3.82%  "
Set-TextValue $ws.Range("D28") "2.28"
$ws.Range("E28").Value = "  -0.77%  "
Set-TextValue $ws.Range("D29") "167.52"
$ws.Range("E29").Value = "  -0.66%  "
Set-TextValue $ws.Range("D30") "34.11"
$ws.Range("E30").Value = "  -0.17%  "
Set-TextValue $ws.Range("D31") "9.18"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +0.13%  "
Set-TextValue $ws.Range("D33") "5.01"
$ws.Range("E33").Value = "  +2.26%  "
Set-TextValue $ws.Range("D34") "4.64"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("E35").Value = "  +4.08%  "
Set-TextValue $ws.Range("D36") "17.06"
$ws.Range("E36").Value = "  +3.20%  "
Set-TextValue $ws.Range("D37") "0.0691"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  +3.81%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D39") "2.82"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D40") "1.79"
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("E41").Value = "  +0.87%  "
Set-TextValue $ws.Range("D42") "2.31"
$ws.Range("E42").Value = "  -6.21%  "
$ws.Range("D43").Value = "2.004.42"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("E44").Value = "  +3.03%  "
Set-TextValue $ws.Range("D45") "10.20"
$ws.Range("E45").Value = "  +6.77%  "
Set-TextValue $ws.Range("D46") "17.47"
$ws.Range("E46").Value = "  +0.23%  "
Set-TextValue $ws.Range("D47") "2.85"
$ws.Range("E47").Value = "  +1.92%  "
Set-TextValue $ws.Range("D48") "55.65"
$ws.Range("E48").Value = "  +6.77%  "
$ws.Range("D49").Value = "2.530.97"
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("E50").Value = "  +4.69%  "
Set-TextValue $ws.Range("D51") "4.57"
$ws.Range("E51").Value = "  +0.92%  "
